$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.981889843940735
$ws.Range("B1").Value = 2.045217990875244
$ws.Range("C1").Value = 7.791342258453369
$ws.Range("D1").Value = 0.9854145646095276
$ws.Range("E1").Value = 0.4802484512329102
